# ---------------------------------------------------------------------------
# Edit: add "PO Forecast" sheet with ds / PO_Forecast / yhat_lower / yhat_upper
# columns, and rename the "Requested quantity" headers on the two existing
# sheets to Weekly_PO_Qty / Monthly_PO_Qty respectively.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# 1) Rename the value-column headers.
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet after the last existing sheet.
$lastIndex  = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header/date formatting from the "Weekly Quantity" sheet so the new
# sheet's styles line up with the workbook's existing style table (bold
# centered header style, and the date-serial number format for column A).
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A62").PasteSpecial(-4122)  # xlPasteFormats

# 3) Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# 4) Data rows (ds serial date, PO_Forecast, yhat_lower, yhat_upper).
$forecastData = @(
    @(44990.99999999999,232,-53.45302935078862,525.7243324336175),
    @(45011.99999999999,232,-55.64996185215574,516.9255678361783),
    @(45046.99999999999,232,-67.95518388922827,527.4858904489338),
    @(45053.99999999999,232,-90.83089924462823,508.4722687230466),
    @(45060.99999999999,232,-61.05217029337783,522.5359536626761),
    @(45067.99999999999,232,-64.28101752920979,512.1536036069256),
    @(45074.99999999999,232,-62.16977976662695,503.2444775437471),
    @(45081.99999999999,232,-88.14254056853845,533.1562449354653),
    @(45088.99999999999,232,-56.70078329489922,539.706659152315),
    @(45095.99999999999,232,-51.53329020624922,516.9978495506302),
    @(45102.99999999999,232,-37.01156815844437,508.5615105461875),
    @(45109.99999999999,232,-54.97817722337592,541.9063548554543),
    @(45116.99999999999,232,-54.75275627227808,520.3513427726003),
    @(45123.99999999999,232,-41.60638551558812,524.9377869416405),
    @(45130.99999999999,232,-59.4231952543119,513.8879374628001),
    @(45137.99999999999,232,-67.40786089392343,511.0496852797079),
    @(45151.99999999999,232,-57.46140285109725,492.6097732299904),
    @(45172.99999999999,232,-50.74077333411766,499.0253388497781),
    @(45179.99999999999,232,-45.70931237572956,501.2314646784566),
    @(45186.99999999999,232,-56.92581564778802,500.7377863692814),
    @(45200.99999999999,232,-62.74373883236316,527.5845815389407),
    @(45207.99999999999,232,-57.88859630263153,514.3240013777429),
    @(45214.99999999999,232,-54.35422275922176,523.0282718487605),
    @(45221.99999999999,232,-23.57204992722915,531.6513041202547),
    @(45228.99999999999,232,-51.91815825850789,531.4567055526849),
    @(45235.99999999999,232,-75.32066505173901,499.9287758162881),
    @(45242.99999999999,232,-48.9494118747411,524.0343932780864),
    @(45249.99999999999,232,-59.9485828744769,513.7605903118889),
    @(45256.99999999999,232,-47.73962474733836,515.3928751106306),
    @(45263.99999999999,232,-59.87246676002444,509.1611303152621),
    @(45277.99999999999,232,-41.60444013415268,519.2363315968707),
    @(45298.99999999999,232,-56.10810592168248,518.8497949997512),
    @(45305.99999999999,232,-72.28118745288744,513.7964571404993),
    @(45312.99999999999,232,-56.02661231695842,516.0102474182471),
    @(45319.99999999999,232,-60.88006956059125,519.6165218729855),
    @(45326.99999999999,232,-40.75900329164643,504.1952653756281),
    @(45382.99999999999,232,-55.51574186933848,518.1264023149578),
    @(45396.99999999999,232,-39.60097157232273,516.292468998497),
    @(45403.99999999999,232,-63.94935300043258,514.7647873798313),
    @(45410.99999999999,232,-46.29390925953145,521.511837982554),
    @(45417.99999999999,232,-64.3137258755003,519.4427980322724),
    @(45424.99999999999,232,-42.41965656920306,512.0898040277323),
    @(45431.99999999999,232,-36.09797506311603,538.8093751885021),
    @(45438.99999999999,232,-56.24648370975486,531.1141338231263),
    @(45459.99999999999,232,-87.55952380071685,532.4776557688997),
    @(45466.99999999999,232,-45.63014252852665,510.5470533777044),
    @(45473.99999999999,232,-49.89607278103605,517.6903463823734),
    @(45487.99999999999,232,-47.81233804635548,522.373006847017),
    @(45522.99999999999,232,-48.03317794480747,526.8592947938666),
    @(45529.99999999999,232,-68.01907759844498,511.9173249367543),
    @(45536.99999999999,232,-77.63142122303671,503.5196692806982),
    @(45543.99999999999,232,-62.21028736517411,502.8250643983228),
    @(45550.99999999999,232,-54.30093787866034,520.7774516014854),
    @(45557.99999999999,232,-70.62790272261765,493.3628454528417),
    @(45564.99999999999,232,-49.43406122919497,521.4098992242938),
    @(45571.99999999999,233,-50.21073563183444,510.493058026571),
    @(45578.99999999999,233,-75.80230712353978,503.3015575000107),
    @(45585.99999999999,233,-45.72236185801946,518.530830542867),
    @(45592.99999999999,233,-37.12560911509621,532.5101307643583),
    @(45599.99999999999,233,-33.17267228662396,500.0310615951012),
    @(45606.99999999999,233,-63.41154100627495,526.481563286412)
)

$r = 2
foreach ($row in $forecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
